$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" column (R) to the essential-services spending table,
# mirroring the formatting already used by the 2020 column (Q).

# R4 - header year 2021 (bold 9pt Times New Roman, medium top+bottom
# border, right/center aligned, same as the other year headers in row 4)
$ws.Range("R4").Value = 2021
$ws.Range("R4").Font.Name = "Times New Roman"
$ws.Range("R4").Font.Size = 9
$ws.Range("R4").Font.Bold = $true
$ws.Range("R4").HorizontalAlignment = -4152
$ws.Range("R4").VerticalAlignment = -4108
$ws.Range("R4").Borders.Item(8).LineStyle = 1
$ws.Range("R4").Borders.Item(8).Weight = -4138
$ws.Range("R4").Borders.Item(9).LineStyle = 1
$ws.Range("R4").Borders.Item(9).Weight = -4138

# R5 - "Education" row value
$ws.Range("R5").Value = 47.8
$ws.Range("R5").Font.Name = "Times New Roman"
$ws.Range("R5").Font.Size = 9
$ws.Range("R5").Font.Bold = $true
$ws.Range("R5").NumberFormat = "0.0"
$ws.Range("R5").HorizontalAlignment = -4152
$ws.Range("R5").VerticalAlignment = -4108

# R6 - "Health" row value
$ws.Range("R6").Value = 20.7
$ws.Range("R6").Font.Name = "Times New Roman"
$ws.Range("R6").Font.Size = 9
$ws.Range("R6").Font.Bold = $false

# R7 - "Social protection" row value
$ws.Range("R7").Value = 9.8
$ws.Range("R7").Font.Name = "Times New Roman"
$ws.Range("R7").Font.Size = 9
$ws.Range("R7").Font.Bold = $false
$ws.Range("R7").NumberFormat = "0.0"

# R8 - bottom total row value, with a medium bottom border closing the
# table
$ws.Range("R8").Value = 17.3
$ws.Range("R8").Font.Name = "Times New Roman"
$ws.Range("R8").Font.Size = 9
$ws.Range("R8").Font.Bold = $false
$ws.Range("R8").NumberFormat = "0.0"
$ws.Range("R8").Borders.Item(9).LineStyle = 1
$ws.Range("R8").Borders.Item(9).Weight = -4138

# Move the active selection the way the author left it (one row below the
# last edited cell, in column P).
$ws.Range("P10").Select() | Out-Null
